# TC_04 - Absence script complete for requesting a single day in current month
#
# Sets the "vacation" sheet's Vacation Start date to June 28, 2022 and
# clears the Vacation End date (single-day absence request), then leaves
# the selection on the now-empty Vacation End cell (H2).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("vacation")

$g2 = $ws.Range("G2")
$h2 = $ws.Range("H2")

# Preserve G2's existing cell formatting (it uses a "quote-prefixed text"
# style, s="8") across the value change by stashing a copy of its format
# in a scratch cell, well outside the sheet's real data, then restoring
# it after the new value is written.
$scratch = $ws.Range("Z100")
$g2.Copy() | Out-Null
$scratch.PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# Update the vacation start date to the new requested day.
$g2.Value = "June 28, 2022"

# Restore G2's original formatting (quote-prefixed text style) that the
# plain value assignment above would otherwise have reset.
$scratch.Copy() | Out-Null
$g2.PasteSpecial(-4122) | Out-Null        # xlPasteFormats
$scratch.Clear() | Out-Null

# A single day absence request only needs the start date populated -
# clear the vacation end date while keeping its existing formatting.
$h2.ClearContents() | Out-Null

# Match the author's final selection: the (now blank) Vacation End cell.
$ws.Activate() | Out-Null
$h2.Select() | Out-Null
